$wb = $excel.ActiveWorkbook

$stdinfo = $wb.Worksheets.Item("StdInfo")
$method1 = $wb.Worksheets.Item("Method1")
$method2 = $wb.Worksheets.Item("Method2")

# --- StdInfo standard concentration/volume corrections ---
# (Method1/Method2 pull these via VLOOKUP, so fixing StdInfo ripples
#  through every dependent row automatically.)

# Row 135: dPA 15:0_18:1
$stdinfo.Range("C135").Value = 1
$stdinfo.Range("D135").Value = 0.1

# Row 137: dCer d18:0/13:0
$stdinfo.Range("C137").Value = 0.01
$stdinfo.Range("D137").Value = 1

# Row 139: dHexCER d18:1/15:0
$stdinfo.Range("C139").Value = 0.2
$stdinfo.Range("D139").Value = 1

# Row 141: dLacCER d18:1/15:0
$stdinfo.Range("C141").Value = 0.1
$stdinfo.Range("D141").Value = 1

# Row 143: dFA 18:1
$stdinfo.Range("C143").Value = 0.1
$stdinfo.Range("D143").Value = 1

# D144's value is unchanged, but it now carries explicit (default) formatting
$stdinfo.Range("D144").NumberFormat = "General"

# Rows 135 and 141 lose their custom 15.75pt height (back to the sheet default)
$stdinfo.Rows("135:135").AutoFit()
$stdinfo.Rows("141:141").AutoFit()

# Small spacer rows newly formatted between the standards (136, 138, 140, 142)
$stdinfo.Range("C136").Font.Bold = $false
$stdinfo.Range("D136").Font.Bold = $false
$stdinfo.Range("C138").Font.Bold = $false
$stdinfo.Range("D138").Font.Bold = $false
$stdinfo.Range("C140").Font.Bold = $false
$stdinfo.Range("D140").Font.Bold = $false
$stdinfo.Range("C142").Font.Bold = $false
$stdinfo.Range("D142").Font.Bold = $false

# --- View-state changes: scroll/selection on each sheet, StdInfo becomes active ---
$method1.Range("E684").Select()
$method2.Range("D528").Select()
$stdinfo.Range("F146").Select()
